$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.44

# Row 3 updates
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 2.75
$ws.Range("N3").Value = 8
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 17
$ws.Range("AC3").Value = 8
$ws.Range("AE3").Value = 17
$ws.Range("AH3").Value = 19
$ws.Range("AZ3").Value = 81

# Row 6 updates
$ws.Range("G6").Value = 3.4
$ws.Range("I6").Value = 2.4
$ws.Range("M6").Value = 1.14
$ws.Range("N6").Value = 5.5
$ws.Range("Z6").Value = 41
$ws.Range("AE6").Value = 21
$ws.Range("BB6").Value = 351

# Row 7 updates
$ws.Range("G7").Value = 2.25
$ws.Range("I7").Value = 3.1
$ws.Range("J7").Value = 3
$ws.Range("L7").Value = 4
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.5
$ws.Range("U7").Value = 1.91
$ws.Range("V7").Value = 1.8
$ws.Range("X7").Value = 10
$ws.Range("AI7").Value = 12
$ws.Range("AK7").Value = 29
$ws.Range("AT7").Value = 2.5
$ws.Range("AX7").Value = 19
$ws.Range("AZ7").Value = 67
$ws.Range("BB7").Value = 251

$wb.Save()
